# Updating some research files
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: replace C29, add D29/E29 (new content, wrapped row) ---
$ws.Range("C29").Value = "Design a Comprehensive Execution Plan for the LLM Cleaner (Phase 2)"
$ws.Range("D29").Value = "While reviewing the existing roadmap, I found it provides a solid high-level direction but lacks actionable detail on how to execute the offline LLM cleaning pipeline—particularly how it integrates with GitHub workflows. Today’s focus was to fill in those missing specifics."
$ws.Range("E29").Value = "A key challenge emerged when attempting to upload the quantized LLM model to GitHub to enable a fully online workflow. GitHub’s file size limits and runner `nconstraints made this infeasible. This broke the assumption that the entire pipeline could run in the cloud autonomously, revealing the need for an `noffline/hybrid execution path that would rely on my device being powered on at certain steps."
$ws.Range("E29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 57.6

# --- Row 30: new date, new task, new status ---
$ws.Range("A30").Value = "1/8/2025(Vacation)"
$ws.Range("C30").Value = "Write the script that will download the model locally"
$ws.Range("F30").Value = "DONE: The code has been written and the model has been downloaded"

# --- Row 31: new date, task carried from old row 32, new status ---
$ws.Range("A31").Value = "1/8/2025(Vacation)"
$ws.Range("C31").Value = "Know the inputs and outputs of each phase and what happens at each phase"
$ws.Range("F31").Value = "DONE: The folder explaining this point is stored at Car_Clinic_Project/docs/Project Understanding Questions.docx"

# --- New rows 38-40: move old rows 29-31 content (tasks 101, 102, 105) under the new date ---
$ws.Range("A38").Value = "1/8/2025(Vacation)"
$ws.Range("B38").Value = "Car Tracking Project"
$ws.Range("C38").Value = "Check how the medical o1 data set looked like so that you can format your data set in the same way"

$ws.Range("A39").Value = "1/8/2025(Vacation)"
$ws.Range("B39").Value = "Car Tracking Project"
$ws.Range("C39").Value = "Improve the file structure and readme file to include all the needed phases."

$ws.Range("A40").Value = "1/8/2025(Vacation)"
$ws.Range("B40").Value = "Car Tracking Project"
$ws.Range("C40").Value = "Generate a fake meta data for the branches so that you can test with later on"

# --- Clear old row 32 (content moved up into rows 30/31) ---
$ws.Range("A32:F32").Clear()

# --- Column F width change (closest value reachable through Excel's pixel-quantised ColumnWidth) ---
$ws.Columns.Item(6).ColumnWidth = 103.3

# --- View/selection updates ---
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C37").Select()
